$d = $word.ActiveDocument

# The "_GoBack" bookmark currently sits right before the "IE8:" run
# (near the end of the compatibility-testing list). It needs to move to
# the very start of the document, right before the "COMPATIBILITY
# TESTING:" run in the first paragraph.
#
# Word's COM model re-anchors a bookmark (rather than erroring) when you
# Bookmarks.Add() a name that already exists elsewhere, so adding it at
# the new location automatically removes it from the old one.
#
# A truly collapsed Range(0,0) at the very start of the document can't
# be used directly as the bookmark anchor here, so a temporary
# placeholder character is inserted at position 0, the bookmark is
# added around that single character, and then the placeholder is
# deleted again - leaving the (now empty) bookmark anchored exactly
# before the first run.
$start = $d.Range(0, 0)
$start.InsertBefore("X")
$anchor = $d.Range(0, 1)
$d.Bookmarks.Add("_GoBack", $anchor)
$d.Range(0, 1).Delete()
